# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# The "Rules" worksheet holds a small lookup table (rows 7-11) mapping an
# hour range to a greeting, with column B holding the rule name (R10, R20,
# R30, R40, ...). This change updates the last rule's name in cell B11
# from "R40" to "1", keeping the value stored as text (as all the other
# rule-name cells in that column are) rather than as a number.
#
# We can't just do $ws.Range("B11").Value = "1" here: Excel's COM layer
# auto-detects numeric-looking strings and stores them as a Number, which
# would change the cell's type. To force a genuine text value without
# disturbing B11's existing cell style/number format, we stage the text
# in an unused scratch cell (explicitly formatted as Text), copy it, and
# paste *values only* into B11 - exactly like using Paste Special > Values
# in the UI - then clean up the scratch cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"

$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)   # xlPasteValues

$scratch.Clear()
